# Updates lat/lon outputs for TT rule (Transitions Rule) summary tables
# Sheet "Means" and "Standard Deviations" - columns D:G (1/3/5/10 mile buffers), rows 2-10

$wb = $excel.ActiveWorkbook

$wsMeans = $wb.Worksheets.Item("Means")
$wsSD    = $wb.Worksheets.Item("Standard Deviations")

# --- Means sheet updates ---
$meansUpdates = @{
  "D2" = 61;   "E2" = 60;   "F2" = 58;   "G2" = 54;
  "D3" = 19;   "E3" = 15;   "F3" = 15;   "G3" = 19;
  "D4" = 20;   "E4" = 25;   "F4" = 26;   "G4" = 27;
  "D5" = 32;   "E5" = 37;   "F5" = 39;   "G5" = 34;
  "D6" = 63;   "E6" = 66;   "F6" = 65;   "G6" = 70;
  "D7" = 8.2;  "E7" = 8.5;  "F7" = 8.5;  "G7" = 7.6;
  "D8" = 8;    "E8" = 7;    "F8" = 6.8;  "G8" = 6.1;
  "D9" = 43;   "E9" = 38;   "F9" = 37;   "G9" = 36;
  "D10" = 0.43; "E10" = 0.43; "F10" = 0.42; "G10" = 0.42;
}

foreach ($addr in $meansUpdates.Keys) {
  $wsMeans.Range($addr).Value = $meansUpdates[$addr]
}

# --- Standard Deviations sheet updates ---
$sdUpdates = @{
  "D2" = 32;   "E2" = 26;   "F2" = 25;   "G2" = 30;
  "D3" = 29;   "E3" = 26;   "F3" = 25;   "G3" = 33;
  "D4" = 19;   "E4" = 19;   "F4" = 19;
  "D5" = 28;   "E5" = 25;
  "D6" = 28;   "E6" = 31;   "F6" = 29;   "G6" = 29;
  "D7" = 9.1;  "E7" = 10;   "F7" = 9.9;  "G7" = 9.4;
  "D8" = 11;   "E8" = 9.4;  "F8" = 8.4;  "G8" = 8.3;
  "D9" = 23;   "E9" = 19;   "F9" = 16;   "G9" = 13;
  "D10" = 0.092; "E10" = 0.088; "F10" = 0.083; "G10" = 0.078;
}

foreach ($addr in $sdUpdates.Keys) {
  $wsSD.Range($addr).Value = $sdUpdates[$addr]
}
